# "correction du titre H1" -- add a new SEO-audit row (row 10) documenting
# the H1 heading-tag issue, matching the structure of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row content -------------------------------------------------
# Column order when first writing each cell matters: it controls the
# order new entries are appended to the shared-string table, and here
# that order is Categorie(A, reused) -> B -> C -> E(reused) -> F -> D.
$ws.Range("A10").Value = "SEO"
$ws.Range("B10").Value = "Titre H1"
$ws.Range("C10").Value = "La balise H1 contient seulement les mots ""La chouette agence."" Il faudrait mettre plus de mots pertinent pour plus d'efficacité."
$ws.Range("E10").Value = "X"
$ws.Range("F10").Value = "smartkeyword - Balise H1"
$ws.Hyperlinks.Add($ws.Range("F10"), "https://smartkeyword.io/optimiser-balise-h1-seo/")
$ws.Range("D10").Value = "Le titre H1 doit contenir des mots clés qui décrivent correctement le site comme par exemple ""La Chouette Agence – Entreprise de webdesign basé à Lyon"""

# --- Formatting: match the look of the rows above it ------------------
$ws.Range("B10:D10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 71

# --- Selection / view state --------------------------------------------
$ws.Range("D10").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Left = 0
